$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sponsor rows added in column H/I (rows 21-24) - write in diff order
$ws.Range("H21").Value = "丨袖手旁观丨"
$ws.Range("I21").Value = 2096927165
$ws.Range("H22").Value = "小牛同志丶"
$ws.Range("I22").Value = 89160614
$ws.Range("H23").Value = "不武酱"
$ws.Range("I23").Value = 416503868
$ws.Range("H24").Value = "蓝瘦香菇mmm"
$ws.Range("I24").Value = 366425370

# New sponsor rows added in column J/K (rows 11-24) - write in diff order
$ws.Range("J11").Value = "以可爱出名#5383"
$ws.Range("K11").Value = 24682425
$ws.Range("J12").Value = "圣灵死法师"
$ws.Range("K12").Value = 838476900
$ws.Range("J13").Value = "SAODIH"
$ws.Range("K13").Value = 21235704
$ws.Range("J14").Value = "玩腻玛个皮皮蛋"
$ws.Range("K14").Value = 259338775
$ws.Range("J15").Value = "玩东方入魔的超#5563"
$ws.Range("K15").Value = 1945313488
$ws.Range("J16").Value = "猪食"
$ws.Range("K16").Value = 185409653
$ws.Range("J17").Value = "小黑爱大白"
$ws.Range("K17").Value = 848895504
$ws.Range("J18").Value = "悲欢离合"
$ws.Range("K18").Value = 970908405
$ws.Range("J19").Value = "nickyly129"
$ws.Range("K19").Value = 1406966725
$ws.Range("J20").Value = "金笙丶"
$ws.Range("K20").Value = 476387019
$ws.Range("J21").Value = "羁绊轮回"
$ws.Range("K21").Value = 1407806903
$ws.Range("J22").Value = "相守不离"
$ws.Range("K22").Value = 39350822
$ws.Range("J23").Value = "郝可怜"
$ws.Range("K23").Value = 947015907
$ws.Range("J24").Value = "亲爱的老王叔"
$ws.Range("K24").Value = 1524326451

# New sponsor rows added in column L/M (rows 11-17) - write in diff order
$ws.Range("L11").Value = "千枫落"
$ws.Range("M11").Value = 1199483482
$ws.Range("L12").Value = "q856101589"
$ws.Range("M12").Value = 85817056
$ws.Range("L13").Value = "XXI1235"
$ws.Range("M13").Value = 1884797690
$ws.Range("L14").Value = "佛丿大湿兄"
$ws.Range("M14").Value = 138245006
$ws.Range("L15").Value = "我是怪叔叔"
$ws.Range("M15").Value = 55883798
$ws.Range("L16").Value = "我才不需要"
$ws.Range("M16").Value = 237209239
$ws.Range("L17").Value = "play"
$ws.Range("M17").Value = 208207478

# Resize the newly populated columns to fit their content
$ws.Columns.Item(10).ColumnWidth = 19.142857142857142
$ws.Columns.Item(11).ColumnWidth = 10.857142857142858
$ws.Columns.Item(12).ColumnWidth = 12.571428571428571
$ws.Columns.Item(13).ColumnWidth = 10.857142857142858

# Update the view state: scroll position and active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 6
$ws.Range("J16").Select() | Out-Null
